$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 -------------------------------------------------------------
# S2 / AX2 hold the order reference number (plain text, already Text-formatted)
$ws.Range("S2").Value = "4192190520"
$ws.Range("AX2").Value = "4192190520"

# T2 holds a date-like label that must stay literal TEXT ("01/08/2016").
# Because T2's cell format is a date format, assigning the string directly
# would be auto-parsed into a real date serial (Excel recognises 01/08/2016
# as a valid M/D/Y date). To avoid that, stage the text in a scratch cell
# that is explicitly Text-formatted, copy it, and paste-special "Values
# only" into T2 - this carries the literal text across without re-parsing
# it and without disturbing T2's existing number format/style.
$scratch = $ws.Range("AG2")
$scratch.NumberFormat = "@"
$scratch.Value = "01/08/2016"
$scratch.Copy()
$ws.Range("T2").PasteSpecial(-4163)
$scratch.Clear()
$excel.CutCopyMode = $false

# --- Row 3 -------------------------------------------------------------
$ws.Range("S3").Value = "4192190520"
$ws.Range("AX3").Value = "4192190520"
$ws.Range("T3").Value = "18/08/2016"
